$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.81602566666667
$ws.Range("H2").Value = 77.44807700000001
$ws.Range("I2").Value = 0.7742517153725241
$ws.Range("J2").Value = 0.7742517153725241
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.806900666666666
$ws.Range("N2").Value = 26.420702
$ws.Range("O2").Value = 0.1733678197953833
$ws.Range("P2").Value = 0.1733678197953834
$ws.Range("Q2").Value = 227.3591736544505
$ws.Range("R2").Value = 2046.232562890054
$ws.Range("S2").Value = 0.1342303318669702
$ws.Range("T2").Value = 0.1342303318669702
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.81602566666667
$ws.Range("H3").Value = 77.44807700000001
$ws.Range("I3").Value = 0.7742517153725241
$ws.Range("J3").Value = 0.7742517153725241
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 18.76689066666667
$ws.Range("N3").Value = 56.30067200000001
$ws.Range("O3").Value = 0.3694347242421866
$ws.Range("P3").Value = 0.3694347242421866
$ws.Range("Q3").Value = 484.4865311341939
$ws.Range("R3").Value = 4360.378780207745
$ws.Range("S3").Value = 0.2860354689626884
$ws.Range("T3").Value = 0.2860354689626884
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.81602566666667
$ws.Range("H4").Value = 77.44807700000001
$ws.Range("I4").Value = 0.7742517153725241
$ws.Range("J4").Value = 0.7742517153725241
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.225144
$ws.Range("N4").Value = 69.675432
$ws.Range("O4").Value = 0.4571974559624301
$ws.Range("P4").Value = 0.4571974559624301
$ws.Range("Q4").Value = 599.5809136160294
$ws.Range("R4").Value = 5396.228222544265
$ws.Range("S4").Value = 0.3539859145428656
$ws.Range("T4").Value = 0.3539859145428656
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.718527666666667
$ws.Range("H5").Value = 8.155583
$ws.Range("I5").Value = 0.08153170965901445
$ws.Range("J5").Value = 0.08153170965901445
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.806900666666666
$ws.Range("N5").Value = 26.420702
$ws.Range("O5").Value = 0.1733678197953833
$ws.Range("P5").Value = 0.1733678197953834
$ws.Range("Q5").Value = 23.94180311991844
$ws.Range("R5").Value = 215.476228079266
$ws.Range("S5").Value = 0.01413497474777353
$ws.Range("T5").Value = 0.01413497474777354
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.718527666666667
$ws.Range("H6").Value = 8.155583
$ws.Range("I6").Value = 0.08153170965901445
$ws.Range("J6").Value = 0.08153170965901445
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 18.76689066666667
$ws.Range("N6").Value = 56.30067200000001
$ws.Range("O6").Value = 0.3694347242421866
$ws.Range("P6").Value = 0.3694347242421866
$ws.Range("Q6").Value = 51.01831149464179
$ws.Range("R6").Value = 459.164803451776
$ws.Range("S6").Value = 0.03012064467487202
$ws.Range("T6").Value = 0.03012064467487202
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.718527666666667
$ws.Range("H7").Value = 8.155583
$ws.Range("I7").Value = 0.08153170965901445
$ws.Range("J7").Value = 0.08153170965901445
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.225144
$ws.Range("N7").Value = 69.675432
$ws.Range("O7").Value = 0.4571974559624301
$ws.Range("P7").Value = 0.4571974559624301
$ws.Range("Q7").Value = 63.13819652631734
$ws.Range("R7").Value = 568.243768736856
$ws.Range("S7").Value = 0.03727609023636889
$ws.Range("T7").Value = 0.03727609023636889
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.808641333333334
$ws.Range("H8").Value = 14.425924
$ws.Range("I8").Value = 0.1442165749684613
$ws.Range("J8").Value = 0.1442165749684613
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.806900666666666
$ws.Range("N8").Value = 26.420702
$ws.Range("O8").Value = 0.1733678197953833
$ws.Range("P8").Value = 0.1733678197953834
$ws.Range("Q8").Value = 42.34922656429422
$ws.Range("R8").Value = 381.143039078648
$ws.Range("S8").Value = 0.02500251318063959
$ws.Range("T8").Value = 0.02500251318063959
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.808641333333334
$ws.Range("H9").Value = 14.425924
$ws.Range("I9").Value = 0.1442165749684613
$ws.Range("J9").Value = 0.1442165749684613
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.76689066666667
$ws.Range("N9").Value = 56.30067200000001
$ws.Range("O9").Value = 0.3694347242421866
$ws.Range("P9").Value = 0.3694347242421866
$ws.Range("Q9").Value = 90.2432461578809
$ws.Range("R9").Value = 812.1892154209282
$ws.Range("S9").Value = 0.05327861060462613
$ws.Range("T9").Value = 0.05327861060462613
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.808641333333334
$ws.Range("H10").Value = 14.425924
$ws.Range("I10").Value = 0.1442165749684613
$ws.Range("J10").Value = 0.1442165749684613
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.225144
$ws.Range("N10").Value = 69.675432
$ws.Range("O10").Value = 0.4571974559624301
$ws.Range("P10").Value = 0.4571974559624301
$ws.Range("Q10").Value = 111.6813874110187
$ws.Range("R10").Value = 1005.132486699168
$ws.Range("S10").Value = 0.06593545118319559
$ws.Range("T10").Value = 0.06593545118319559